$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Issue 11 (row 35): mark resolved and record close date ---
$ws.Range("B35").Value = "RESOLVED"
$ws.Range("F35").Value = 45494

# --- Augment the existing last history row (row 45) with a new update entry ---
$ws.Range("G45").Value = 45482
$ws.Range("I45").Value = "OPEN"

# --- Append two new history rows at the bottom (before inserting the new row above,
#     so that the new shared strings get appended to the shared-string table in the
#     same order the rows will finally appear on screen: 46 -> 47) ---
$ws.Range("G46").Value = 45483
$ws.Range("H46").Value = "Updated problem diagram and argument based on Sean's feedback"
$ws.Range("I46").Value = "OPEN"
$ws.Rows(46).RowHeight = 17

$ws.Range("G47").Value = 45494
$ws.Range("H47").Value = "Rewrite problem diagram and argument based on Sean's feedback"
$ws.Range("I47").Value = "OPEN"
$ws.Rows(47).RowHeight = 17

# --- Insert a brand-new row right after row 35 for the "process argument eliminated"
#     history entry; this pushes everything from row 36 down by one. ---
$ws.Rows("36:36").Insert()
$ws.Range("A36").Clear()
$ws.Range("C36").Clear()
$ws.Range("E36").Clear()
$ws.Range("F36").Clear()
$ws.Range("G36").Value = 45494
$ws.Range("H36").Value = "Process argument eliminated from problem argument in favor of direct reasoning"
$ws.Range("I36").Value = "RESOLVED"
$ws.Rows(36).RowHeight = 17

# --- Update the view selection to reflect where the editor last left off ---
$ws.Range("I37").Select()
